$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.538.64"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "2.579.09"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "542.12"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "144.06"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.583"
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("D9").Value = "6.76"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  -2.99%  "
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "3.032.08"
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("D14").Value = "58.454.62"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "20.55"
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000131"
$ws.Range("E16").Value = "  -2.77%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.542.11"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").Value = "4.46"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "333.67"
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").Value = "10.04"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "6.09"
$ws.Range("E21").Value = "  -4.02%  "
$ws.Range("D23").Value = "66.36"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "0.422"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -4.91%  "
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("D28").Value = "0.0₃0739"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "6.00"
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("D32").Value = "152.94"
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("D33").Value = "18.92"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "3.91"
$ws.Range("E34").Value = "  -3.00%  "
$ws.Range("D35").Value = "0.849"
$ws.Range("E35").Value = "  +2.86%  "
$ws.Range("E36").Value = "  -4.82%  "
$ws.Range("D37").Value = "0.821"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("D38").Value = "1.42"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").Value = "278.45"
$ws.Range("E40").Value = "  -5.00%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  -2.32%  "
$ws.Range("D43").Value = "10.63"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").Value = "0.0943"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").Value = "18.50"
$ws.Range("E46").Value = "  -5.12%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "1.900.60"
$ws.Range("E48").Value = "  -3.74%  "
$ws.Range("D49").Value = "17.85"
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("D51").Value = "109.05"
$ws.Range("E51").Value = "  -1.24%  "
